$wb = $excel.ActiveWorkbook

# --- Remove header-row styling (bold/border/center) on all sheets ---
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("A1:N1").ClearFormats()
}

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 1566.6666
$ws.Cells.Item(2, 9).Value = 1566.6666
$ws.Cells.Item(2, 11).Value = 1566.6666
$ws.Cells.Item(2, 13).Value = -1453.6666
$ws.Cells.Item(43, 8).Value = 1705.5834
$ws.Cells.Item(43, 9).Value = 1395
$ws.Cells.Item(43, 10).Value = 1809.1111
$ws.Cells.Item(43, 11).Value = 1395
$ws.Cells.Item(43, 12).Value = 1809.1111
$ws.Cells.Item(43, 13).Value = -1326
$ws.Cells.Item(43, 14).Value = -1947.1111
$ws.Cells.Item(70, 8).Value = 5717907.5
$ws.Cells.Item(70, 10).Value = 2304.8635
$ws.Cells.Item(70, 12).Value = 6914.5905
$ws.Cells.Item(70, 14).Value = -7454.5905
$ws.Cells.Item(73, 8).Value = 5717907.5
$ws.Cells.Item(73, 10).Value = 2304.8635
$ws.Cells.Item(73, 12).Value = 6914.5905
$ws.Cells.Item(73, 14).Value = -8786.5905
$ws.Cells.Item(74, 8).Value = 7817
$ws.Cells.Item(74, 9).Value = 7225.5
$ws.Cells.Item(74, 11).Value = 7225.5
$ws.Cells.Item(74, 13).Value = -6289.5
$ws.Cells.Item(77, 8).Value = 7817
$ws.Cells.Item(77, 9).Value = 7225.5
$ws.Cells.Item(77, 11).Value = 36127.5
$ws.Cells.Item(77, 13).Value = -31447.5
$ws.Cells.Item(116, 8).Value = 3899.5
$ws.Cells.Item(116, 10).Value = 4279.4
$ws.Cells.Item(116, 12).Value = 4279.4
$ws.Cells.Item(116, 14).Value = -11163.4
$ws.Cells.Item(132, 8).Value = 4217.0356
$ws.Cells.Item(132, 9).Value = 4500.5
$ws.Cells.Item(132, 11).Value = 13501.5
$ws.Cells.Item(132, 13).Value = -10971.5
$ws.Cells.Item(134, 8).Value = 86580.57000000001
$ws.Cells.Item(134, 10).Value = 86580.57000000001
$ws.Cells.Item(134, 12).Value = 86580.57000000001
$ws.Cells.Item(134, 14).Value = -96720.57000000001
$ws.Cells.Item(137, 8).Value = 5016941
$ws.Cells.Item(137, 9).Value = 1031896
$ws.Cells.Item(137, 11).Value = 3095688
$ws.Cells.Item(137, 13).Value = -3093138

# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 45587024
$ws.Cells.Item(74, 9).Value = 62680404
$ws.Cells.Item(74, 10).Value = 4671.3335
$ws.Cells.Item(74, 11).Value = 62680404
$ws.Cells.Item(74, 12).Value = 4671.3335
$ws.Cells.Item(74, 13).Value = -62679530
$ws.Cells.Item(74, 14).Value = -6419.3335
$ws.Cells.Item(77, 8).Value = 45587024
$ws.Cells.Item(77, 9).Value = 62680404
$ws.Cells.Item(77, 10).Value = 4671.3335
$ws.Cells.Item(77, 11).Value = 313402020
$ws.Cells.Item(77, 12).Value = 23356.6675
$ws.Cells.Item(77, 13).Value = -313397652
$ws.Cells.Item(77, 14).Value = -32092.6675
$ws.Cells.Item(92, 8).Value = 137528860
$ws.Cells.Item(92, 10).Value = 137528860
$ws.Cells.Item(92, 12).Value = 137528860
$ws.Cells.Item(92, 14).Value = -137533852
$ws.Cells.Item(97, 8).Value = 7375.294
$ws.Cells.Item(97, 9).Value = 8292.429
$ws.Cells.Item(97, 11).Value = 8292.429
$ws.Cells.Item(97, 13).Value = -7796.429
$ws.Cells.Item(122, 8).Value = 741426.5600000001
$ws.Cells.Item(122, 9).Value = 3578
$ws.Cells.Item(122, 10).Value = 2807402.5
$ws.Cells.Item(122, 11).Value = 10734
$ws.Cells.Item(122, 12).Value = 8422207.5
$ws.Cells.Item(122, 13).Value = -8284
$ws.Cells.Item(122, 14).Value = -8427107.5
$ws.Cells.Item(132, 8).Value = 3210.1482
$ws.Cells.Item(132, 9).Value = 2920.739
$ws.Cells.Item(132, 11).Value = 8762.217000000001
$ws.Cells.Item(132, 13).Value = -6232.217000000001

# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 4391.4165
$ws.Cells.Item(20, 9).Value = 1843.25
$ws.Cells.Item(20, 11).Value = 1843.25
$ws.Cells.Item(20, 13).Value = -1596.25
$ws.Cells.Item(105, 8).Value = 113497.6
$ws.Cells.Item(105, 9).Value = 159139.58
$ws.Cells.Item(105, 11).Value = 159139.58
$ws.Cells.Item(105, 13).Value = -157392.58
$ws.Cells.Item(141, 8).Value = 109500.5
$ws.Cells.Item(141, 10).Value = 109500.5
$ws.Cells.Item(141, 12).Value = 109500.5
$ws.Cells.Item(141, 14).Value = -119860.5

# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2248.7737
$ws.Cells.Item(31, 9).Value = 730.7857
$ws.Cells.Item(31, 10).Value = 2793.6924
$ws.Cells.Item(31, 11).Value = 730.7857
$ws.Cells.Item(31, 12).Value = 2793.6924
$ws.Cells.Item(31, 13).Value = -435.7857
$ws.Cells.Item(31, 14).Value = -3383.6924
$ws.Cells.Item(34, 8).Value = 2248.7737
$ws.Cells.Item(34, 9).Value = 730.7857
$ws.Cells.Item(34, 10).Value = 2793.6924
$ws.Cells.Item(34, 11).Value = 730.7857
$ws.Cells.Item(34, 12).Value = 2793.6924
$ws.Cells.Item(34, 13).Value = -528.7857
$ws.Cells.Item(34, 14).Value = -3197.6924
$ws.Cells.Item(105, 8).Value = 163850.84
$ws.Cells.Item(105, 9).Value = 211907.1
$ws.Cells.Item(105, 10).Value = 3663.3333
$ws.Cells.Item(105, 11).Value = 211907.1
$ws.Cells.Item(105, 12).Value = 3663.3333
$ws.Cells.Item(105, 13).Value = -210160.1
$ws.Cells.Item(105, 14).Value = -7157.3333
$ws.Cells.Item(107, 8).Value = 52643920
$ws.Cells.Item(107, 9).Value = 66681584
$ws.Cells.Item(107, 10).Value = 2675
$ws.Cells.Item(107, 11).Value = 66681584
$ws.Cells.Item(107, 12).Value = 2675
$ws.Cells.Item(107, 13).Value = -66679664
$ws.Cells.Item(107, 14).Value = -6515
$ws.Cells.Item(134, 8).Value = 3298267.5
$ws.Cells.Item(134, 9).Value = 4818738.5
$ws.Cells.Item(134, 11).Value = 14456215.5
$ws.Cells.Item(134, 13).Value = -14453680.5

# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(7, 8).Value = 111389130
$ws.Cells.Item(7, 9).Value = 416953.34
$ws.Cells.Item(7, 11).Value = 1250860.02
$ws.Cells.Item(7, 13).Value = -1250748.02
$ws.Cells.Item(98, 8).Value = 745.2941
$ws.Cells.Item(98, 10).Value = 613.8182
$ws.Cells.Item(98, 12).Value = 1841.4546
$ws.Cells.Item(98, 14).Value = -4837.4546
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).ClearContents()

# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 6942.579
$ws.Cells.Item(70, 9).Value = 7266.8667
$ws.Cells.Item(70, 10).Value = 5726.5
$ws.Cells.Item(70, 11).Value = 7266.8667
$ws.Cells.Item(70, 12).Value = 5726.5
$ws.Cells.Item(70, 13).Value = -6996.8667
$ws.Cells.Item(70, 14).Value = -6266.5
$ws.Cells.Item(73, 8).Value = 6942.579
$ws.Cells.Item(73, 9).Value = 7266.8667
$ws.Cells.Item(73, 10).Value = 5726.5
$ws.Cells.Item(73, 11).Value = 7266.8667
$ws.Cells.Item(73, 12).Value = 5726.5
$ws.Cells.Item(73, 13).Value = -6330.8667
$ws.Cells.Item(73, 14).Value = -7598.5
$ws.Cells.Item(80, 8).Value = 5720.7617
$ws.Cells.Item(80, 10).Value = 3626.923
$ws.Cells.Item(80, 12).Value = 3626.923
$ws.Cells.Item(80, 14).Value = -5622.923
$ws.Cells.Item(83, 8).Value = 5720.7617
$ws.Cells.Item(83, 10).Value = 3626.923
$ws.Cells.Item(83, 12).Value = 18134.615
$ws.Cells.Item(83, 14).Value = -28118.615
$ws.Cells.Item(93, 8).Value = 38405
$ws.Cells.Item(93, 10).Value = 38405
$ws.Cells.Item(93, 12).Value = 38405
$ws.Cells.Item(93, 14).Value = -42149
$ws.Cells.Item(99, 8).Value = 9560
$ws.Cells.Item(99, 9).Value = 6566.6665
$ws.Cells.Item(99, 11).Value = 6566.6665
$ws.Cells.Item(99, 13).Value = -4320.6665
$ws.Cells.Item(102, 8).Value = 7468.2
$ws.Cells.Item(102, 9).Value = 7766.875
$ws.Cells.Item(102, 11).Value = 7766.875
$ws.Cells.Item(102, 13).Value = -6144.875
$ws.Cells.Item(107, 8).Value = 473.92856
$ws.Cells.Item(107, 9).Value = 516.5
$ws.Cells.Item(107, 10).Value = 442
$ws.Cells.Item(107, 11).Value = 516.5
$ws.Cells.Item(107, 12).Value = 442
$ws.Cells.Item(107, 13).Value = 1403.5
$ws.Cells.Item(107, 14).Value = -4282
$ws.Cells.Item(111, 8).Value = 11100
$ws.Cells.Item(111, 10).Value = 11100
$ws.Cells.Item(111, 12).Value = 11100
$ws.Cells.Item(126, 8).Value = 8583
$ws.Cells.Item(126, 9).Value = 9349.75
$ws.Cells.Item(126, 10).Value = 7614.4736
$ws.Cells.Item(126, 11).Value = 28049.25
$ws.Cells.Item(126, 12).Value = 22843.4208
$ws.Cells.Item(126, 13).Value = -25579.25
$ws.Cells.Item(126, 14).Value = -27783.4208
$ws.Cells.Item(111, 14).Value = -17234

# --- Sheet 7 (LTW) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(46, 8).Value = 2100.1365
$ws.Cells.Item(46, 9).Value = 890.375
$ws.Cells.Item(46, 10).Value = 2791.4285
$ws.Cells.Item(46, 11).Value = 890.375
$ws.Cells.Item(46, 12).Value = 2791.4285
$ws.Cells.Item(46, 13).Value = -702.375
$ws.Cells.Item(46, 14).Value = -3167.4285
$ws.Cells.Item(50, 8).Value = 50035.5
$ws.Cells.Item(50, 10).Value = 69995
$ws.Cells.Item(50, 12).Value = 69995
$ws.Cells.Item(50, 14).Value = -71269
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(68, 8).Value = 4115.8237
$ws.Cells.Item(68, 9).Value = 2744.8
$ws.Cells.Item(68, 11).Value = 2744.8
$ws.Cells.Item(68, 13).Value = -1995.8
$ws.Cells.Item(71, 8).Value = 4115.8237
$ws.Cells.Item(71, 9).Value = 2744.8
$ws.Cells.Item(71, 11).Value = 13724
$ws.Cells.Item(71, 13).Value = -9980
$ws.Cells.Item(82, 8).Value = 2893.3
$ws.Cells.Item(82, 10).Value = 2294.2222
$ws.Cells.Item(82, 12).Value = 2294.2222
$ws.Cells.Item(82, 14).Value = -3016.2222
$ws.Cells.Item(85, 8).Value = 2893.3
$ws.Cells.Item(85, 10).Value = 2294.2222
$ws.Cells.Item(85, 12).Value = 2294.2222
$ws.Cells.Item(85, 14).Value = -4790.2222
$ws.Cells.Item(116, 8).Value = 256091
$ws.Cells.Item(116, 10).Value = 256091
$ws.Cells.Item(116, 12).Value = 256091
$ws.Cells.Item(116, 14).Value = -265269
$ws.Cells.Item(122, 8).Value = 4476.7334
$ws.Cells.Item(122, 9).Value = 3942.3845
$ws.Cells.Item(122, 11).Value = 11827.1535
$ws.Cells.Item(122, 13).Value = -9377.1535
$ws.Cells.Item(132, 8).Value = 17169.54
$ws.Cells.Item(132, 9).Value = 46412.25
$ws.Cells.Item(132, 10).Value = 4172.778
$ws.Cells.Item(132, 11).Value = 139236.75
$ws.Cells.Item(132, 12).Value = 12518.334
$ws.Cells.Item(132, 13).Value = -136706.75
$ws.Cells.Item(132, 14).Value = -17578.334
$ws.Cells.Item(51, 14).ClearContents()

# --- Sheet 8 (WVR) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(41, 8).Value = 12036
$ws.Cells.Item(41, 10).Value = 11290.5
$ws.Cells.Item(41, 12).Value = 11290.5
$ws.Cells.Item(41, 14).Value = -12070.5
$ws.Cells.Item(81, 8).Value = 11266.7
$ws.Cells.Item(81, 9).Value = 11517.579
$ws.Cells.Item(81, 11).Value = 23035.158
$ws.Cells.Item(81, 13).Value = -21974.158
$ws.Cells.Item(84, 8).Value = 11266.7
$ws.Cells.Item(84, 9).Value = 11517.579
$ws.Cells.Item(84, 11).Value = 115175.79
$ws.Cells.Item(84, 13).Value = -109871.79
$ws.Cells.Item(107, 8).Value = 2384.5557
$ws.Cells.Item(107, 9).Value = 2172.5
$ws.Cells.Item(107, 10).Value = 2808.6667
$ws.Cells.Item(107, 11).Value = 6517.5
$ws.Cells.Item(107, 12).Value = 8426.000100000001
$ws.Cells.Item(107, 13).Value = -4597.5
$ws.Cells.Item(107, 14).Value = -12266.0001
$ws.Cells.Item(122, 8).Value = 5579.3076
$ws.Cells.Item(122, 9).Value = 1941
$ws.Cells.Item(122, 11).Value = 5823
$ws.Cells.Item(122, 13).Value = -3373
$ws.Cells.Item(132, 8).Value = 34035.5
$ws.Cells.Item(132, 9).Value = 43365.785
$ws.Cells.Item(132, 11).Value = 130097.355
$ws.Cells.Item(132, 13).Value = -127567.355
$ws.Cells.Item(140, 8).Value = 81899.39999999999
$ws.Cells.Item(140, 10).Value = 81899.39999999999
$ws.Cells.Item(140, 12).Value = 81899.39999999999
$ws.Cells.Item(140, 14).Value = -92259.39999999999
